$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at B (old B..K shift to C..L), carrying column
#    A's width/style - because the insertion point sits *inside* the
#    existing "uniqueness data" band (cols A:B, style "Bad"/red), the
#    band grows to A:C in place, and the A4:B4 merge auto-extends to
#    A4:C4 too - both exactly matching the target layout already.
$ws.Columns.Item(1).Copy()
$ws.Columns.Item(2).Insert()

# 2. Rows 1-3 are single full-width merged banners; the column insert
#    already auto-extended their merges (A1:K1->A1:L1, etc.) but the
#    newly inserted column B carries a stray copy of column A's text, so
#    clear the whole merged block (exact extent - no Merge() call needed)
#    before writing the banner text back into the anchor cell.
$ws.Range("A1:L1").ClearContents()
$ws.Range("A1").Value = "特征清单"
$ws.Range("A2:L2").ClearContents()
$ws.Range("A2").Value = "为每个中类分别训练一个学习器"
$ws.Range("A3:L3").ClearContents()
$ws.Range("A3").Value = "将每一个中类在每一天的销售数据作为一个特征"

# 3. Row 4 group headers. A4:C4 ("唯一性数据") and J4:K4 ("全局类") kept
#    their correct shape through the insert, so just clear + rewrite
#    them directly. The middle two groups changed shape
#    (C4:E4->D4:F4 needs to become D4:G4; F4:H4->G4:I4 needs to become
#    H4:I4) so unmerge then re-merge those two to their new extent.
$ws.Range("A4:C4").ClearContents()
$ws.Range("A4").Value = "唯一性数据"
$ws.Range("J4:K4").ClearContents()
$ws.Range("J4").Value = "全局类"
$ws.Range("L4").Value = "标签"

$ws.Range("D4:I4").UnMerge()
$ws.Range("D4:G4").Merge()
$ws.Range("D4").Value = "时间类"
$ws.Range("H4:I4").Merge()
$ws.Range("H4").Value = "促销类"

# 4. Row 5 - individual column headers
$ws.Range("A5").Value = "中类编号"
$ws.Range("B5").Value = "总日期"
$ws.Range("C5").Value = "月份"
$ws.Range("D5").Value = "日期"
$ws.Range("E5").Value = "星期"
$ws.Range("F5").Value = "节假日"
$ws.Range("G5").Value = "节假日前"
$ws.Range("H5").Value = "是否有促销活动"
$ws.Range("I5").Value = "大类中其他促销的中类数量"
$ws.Range("J5").Value = "当天总销量"
$ws.Range("K5").Value = "当天总营业额"
$ws.Range("L5").Value = "当天销量"

# 5. Match the saved selection/active cell of the edited file.
$ws.Range("D6:I6").Select()
